$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.892.55'
$ws.Range("E2").Value = '  +0.82%  '

$ws.Range("D3").Value = '3.364.28'
$ws.Range("E3").Value = '  +0.73%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '554.75'
$ws.Range("E5").Value = '  +0.39%  '

$ws.Range("D6").Value = '173.81'
$ws.Range("E6").Value = '  -0.95%  '

$ws.Range("D7").Value = '0.632'
$ws.Range("E7").Value = '  +2.27%  '

$ws.Range("D8").Value = '3.354.61'
$ws.Range("E8").Value = '  +0.64%  '

$ws.Range("E9").Value = '  +0.14%  '

$ws.Range("E10").Value = '  +6.91%  '

$ws.Range("D11").Value = '0.636'
$ws.Range("E11").Value = '  +1.44%  '

$ws.Range("D12").Value = '53.50'
$ws.Range("E12").Value = '  -2.18%  '

$ws.Range("E13").Value = '  +3.16%  '

$ws.Range("E14").Value = '  +0.97%  '

$ws.Range("D15").Value = '3.904.03'
$ws.Range("E15").Value = '  +0.86%  '

$ws.Range("E16").Value = '  +2.24%  '

$ws.Range("D17").Value = '18.21'
$ws.Range("E17").Value = '  -0.58%  '

$ws.Range("D18").Value = '3.364.01'
$ws.Range("E18").Value = '  +0.84%  '

$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '11.86'
$ws.Range("E19").Value = '  +0.99%  '

$ws.Range("B20").Value = 'WrappedBTC'
$ws.Range("C20").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D20").Value = '64.787.58'
$ws.Range("E20").Value = '  +0.84%  '

$ws.Range("E21").Value = '  +1.65%  '

$ws.Range("D22").Value = '451.19'
$ws.Range("E22").Value = '  +4.18%  '

$ws.Range("D23").Value = '4.91'
$ws.Range("E23").Value = '  -4.19%  '

$ws.Range("D24").Value = '4.06'
$ws.Range("E24").Value = '  +0.13%  '

$ws.Range("D25").Value = '87.04'
$ws.Range("E25").Value = '  +3.24%  '

$ws.Range("D26").Value = '13.64'
$ws.Range("E26").Value = '  +1.31%  '

$ws.Range("D27").Value = '10.71'
$ws.Range("E27").Value = '  -0.31%  '

$ws.Range("D28").Value = '2.86'
$ws.Range("E28").Value = '  +0.94%  '

$ws.Range("E29").Value = '  -1.17%  '

$ws.Range("D30").Value = '31.03'
$ws.Range("E30").Value = '  +4.38%  '

$ws.Range("D31").Value = '6.54'
$ws.Range("E31").Value = '  -1.91%  '

$ws.Range("D32").Value = '63.05'
$ws.Range("E32").Value = '  +8.24%  '

$ws.Range("D33").Value = '11.42'
$ws.Range("E33").Value = '  -0.55%  '

$ws.Range("D34").Value = '575.42'
$ws.Range("E34").Value = '  -0.70%  '

$ws.Range("E35").Value = '  -0.31%  '

$ws.Range("E36").Value = '  -0.07%  '

$ws.Range("D37").Value = '3.62'
$ws.Range("E37").Value = '  +3.90%  '

$ws.Range("D38").Value = '0.141'
$ws.Range("E38").Value = '  -0.13%  '

$ws.Range("D39").Value = '35.58'
$ws.Range("E39").Value = '  -0.20%  '

$ws.Range("E40").Value = '  +0.83%  '

$ws.Range("E41").Value = '  -1.66%  '

$ws.Range("D42").Value = '3.084.58'
$ws.Range("E42").Value = '  -0.79%  '

$ws.Range("D43").Value = '0.0416'
$ws.Range("E43").Value = '  +1.97%  '

$ws.Range("E44").Value = '  -1.21%  '

$ws.Range("E45").Value = '  +2.99%  '

$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '3.17'
$ws.Range("E46").Value = '  -1.07%  '

$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").Value = '2.44'
$ws.Range("E47").Value = '  -0.79%  '

$ws.Range("D48").Value = '142.42'
$ws.Range("E48").Value = '  +5.54%  '

$ws.Range("E49").Value = '  +0.16%  '

$ws.Range("E50").Value = '  -2.43%  '

$ws.Range("D51").Value = '8.26'
$ws.Range("E51").Value = '  -0.10%  '
